# Update "想去人数" (want-to-go count) values in the "展览" (sheet1) and
# "全部类型" (sheet4) sheets to match the refreshed export data.

$wb = $excel.ActiveWorkbook

# Row -> new F (想去人数) value, as they appear on the "展览" sheet.
$exhibitionUpdates = @{
    2  = 3
    3  = 12673
    5  = 76
    6  = 45
    9  = 1
    10 = 12567
    11 = 251
    12 = 3
    13 = 4933
    14 = 4876
    15 = 167
    16 = 78
    18 = 114
    19 = 973
    22 = 370
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $ws1.Range("F$row").Value = $exhibitionUpdates[$row]
}

# "全部类型" combines all event types; the same events sit one row lower
# (row 2 there is a "演出" entry that precedes the "展览" rows).
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $exhibitionUpdates.Keys) {
    $targetRow = $row + 1
    $ws4.Range("F$targetRow").Value = $exhibitionUpdates[$row]
}
